$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260. This shifts the existing rows 260-317
# down to 261-318 (Excel copies formatting, e.g. the date style on column D,
# from the row above automatically), and extends the used range/dimension
# to A1:R318.
$ws.Rows("260:260").Insert()

# Populate the newly inserted row 260 with the new weekly price record.
$ws.Range("A260").Value = 4
$ws.Range("B260").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C260").Value = "Los Lagos"
$ws.Range("D260").Value = 44785
$ws.Range("E260").Value = 10
$ws.Range("F260").Value = 100112037
$ws.Range("G260").Value = "Cebollín"
$ws.Range("H260").Value = "Sin especificar"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 200
$ws.Range("K260").Value = 9000
$ws.Range("L260").Value = 10000
$ws.Range("M260").Value = 9500
$ws.Range("N260").Value = "$/paquete 36 unidades"
$ws.Range("O260").Value = "Región Metropolitana"
$ws.Range("P260").Value = 264
$ws.Range("Q260").Value = 36
$ws.Range("R260").Value = "Hortaliza"
